$d = $word.ActiveDocument

# Locate the paragraph holding the concatenated highlight-color string.
$rng = $d.Content
$found = $rng.Find.Execute("#c885da#f9cd59#7cc867#fb5b89#69aff0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $startIndex = $rng.Paragraphs.First.Index

    $values = @("#c885da: 29", "#f9cd59: 20", "#7cc867: 65", "#fb5b89: 50", "#69aff0: 7")

    for ($i = 0; $i -lt $values.Length; $i++) {
        $p = $d.Paragraphs.Item($startIndex + $i)
        $pr = $p.Range
        $pr.Text = $values[$i]
        if ($i -lt ($values.Length - 1)) {
            $pr.InsertParagraphAfter()
        }
    }
}
